# GreenField Properties Proposal - "complete security measures headline"
#
# 1. Trim the dangling "The plan is to..." from the cover-letter paragraph.
# 2. Fill in the Security Measures section's placeholder with real copy.
# 3. Clear out every other remaining "<Enter ... here>" placeholder so the
#    paragraph keeps an empty run (run survives, just loses its <w:t>).

$d = $word.ActiveDocument

# Namespace fragment used to splice an empty run into a paragraph while
# preserving the surrounding <w:p>/<w:pPr> (InsertXML replaces the range's
# contents in place, it does not insert a new paragraph).
$emptyRunXml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr/></w:r></w:p></pkg:xmlData>'

# --- 1. Trim the introductory paragraph -----------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*layed out in your new building. The plan is to...*") {
        $p.Range.Text = "I have enclosed in this document my initial recommendations for how the network should be layed out in your new building. "
    }
}

# --- 2. Write the real Security Measures copy ------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*<Enter your Network Security Proposal here>*") {
        $p.Range.Text = "Greenfield is a business that relies its network infrastructure for its daily operations. It is of upmost importance to fortify and secure the network against unauthorized access (both inside and outside the organization), malware and data breaches. I have below some recommendations to ensure a secure and optimal running network."
    }
}

# --- 3. Blank out the remaining placeholders, keeping an empty run --------
$placeholders = @(
    "<Enter an introduction to your proposal here>",
    "<Enter your Network Infrastructure Proposal here>",
    "<Enter your Network Segmentation Proposal here>",
    "<Enter your Network Wi-Fi Proposal here>",
    "<Enter your Network Physical Security Proposal here>",
    "<Enter your Network Authentication Proposal here>",
    "<Enter your Network Lockout Policy Proposal here>",
    "<Enter your Network Password Complexity Requirements Proposal here>",
    "<Enter your Network Firewall Proposal here>",
    "<Enter your Network Anti-Malware Proposal here>"
)

foreach ($needle in $placeholders) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($needle + "`r")) {
            $r = $d.Range($p.Range.Start, $p.Range.End - 1)
            $r.InsertXML($emptyRunXml)
        }
    }
}
